# Contest 32 PBKS vs RR.
# Fill in the match results for Contest 32 (row 41 on Sheet1).
# Columns E, H, K, N, Q, T hold each team's raw score for the match;
# columns D, G, J, M, P, S already contain VLOOKUP formulas that
# translate the rank of each score into a point value, and will
# recalculate automatically once the raw scores below are entered.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E41").Value = 60
$ws.Range("H41").Value = 80
$ws.Range("K41").Value = 40
$ws.Range("N41").Value = 20
$ws.Range("Q41").Value = 0
$ws.Range("T41").Value = 100

$excel.Calculate()
